$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.921.79'
$ws.Range("D2").Style = "Normal"

$ws.Range("E2").Value = '  -2.25%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.419.88'
$ws.Range("D3").Style = "Normal"

$ws.Range("E3").Value = '  -1.46%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.997'
$ws.Range("D4").Style = "Normal"

$ws.Range("E4").Value = '  -0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '570.94'
$ws.Range("D5").Style = "Normal"

$ws.Range("E5").Value = '  -0.71%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '139.97'
$ws.Range("D6").Style = "Normal"

$ws.Range("E6").Value = '  -2.95%  '

$ws.Range("E7").Value = '  +0.24%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.529'
$ws.Range("D8").Style = "Normal"

$ws.Range("E8").Value = '  -0.74%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.404.84'
$ws.Range("D9").Style = "Normal"

$ws.Range("E9").Value = '  -1.90%  '

$ws.Range("E10").Value = '  -1.01%  '

$ws.Range("E11").Value = '  +0.01%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.10'
$ws.Range("D12").Style = "Normal"

$ws.Range("E12").Value = '  -2.26%  '

$ws.Range("E13").Value = '  -0.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '25.96'
$ws.Range("D14").Style = "Normal"

$ws.Range("E14").Value = '  -1.45%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.0000172'
$ws.Range("D15").Style = "Normal"

$ws.Range("E15").Value = '  -1.68%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '2.818.72'
$ws.Range("D16").Style = "Normal"

$ws.Range("E16").Value = '  -2.67%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '60.868.22'
$ws.Range("D17").Style = "Normal"

$ws.Range("E17").Value = '  -1.93%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.398.54'
$ws.Range("D18").Style = "Normal"

$ws.Range("E18").Value = '  -2.27%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '10.62'
$ws.Range("D19").Style = "Normal"

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '7.41'
$ws.Range("D20").Style = "Normal"

$ws.Range("E20").Value = '  +3.65%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '322.91'
$ws.Range("D21").Style = "Normal"

$ws.Range("E21").Value = '  -1.69%  '

$ws.Range("E22").Value = '  -1.35%  '

$ws.Range("E23").Value = '  +1.12%  '

$ws.Range("E24").Value = '  +0.26%  '

$ws.Range("E25").Value = '  -3.41%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '65.02'
$ws.Range("D26").Style = "Normal"

$ws.Range("E26").Value = '  -1.12%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.52'
$ws.Range("D27").Style = "Normal"

$ws.Range("E27").Value = '  -7.51%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '577.09'
$ws.Range("D28").Style = "Normal"

$ws.Range("E28").Value = '  -2.97%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.515.86'
$ws.Range("D29").Style = "Normal"

$ws.Range("E29").Value = '  -2.29%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0922'
$ws.Range("D30").Style = "Normal"

$ws.Range("E30").Value = '  -3.90%  '

$ws.Range("E31").Value = '  -0.99%  '

$ws.Range("E32").Value = '  -5.51%  '

$ws.Range("E33").Value = '  -1.82%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.133'
$ws.Range("D34").Style = "Normal"

$ws.Range("E34").Value = '  -3.04%  '

$ws.Range("E35").Value = '  +0.09%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '4.66'
$ws.Range("D36").Style = "Normal"

$ws.Range("E36").Value = '  -5.08%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.41'
$ws.Range("D37").Style = "Normal"

$ws.Range("E37").Value = '  -2.80%  '

$ws.Range("E38").Value = '  -2.38%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '149.18'
$ws.Range("D39").Style = "Normal"

$ws.Range("E39").Value = '  -1.55%  '

$ws.Range("E40").Value = '  -0.82%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.15'
$ws.Range("D41").Style = "Normal"

$ws.Range("E41").Value = '  -3.94%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.00'
$ws.Range("D42").Style = "Normal"

$ws.Range("E43").Value = '  -3.51%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '41.07'
$ws.Range("D44").Style = "Normal"

$ws.Range("E44").Value = '  -3.58%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.33'
$ws.Range("D45").Style = "Normal"

$ws.Range("E45").Value = '  -6.06%  '

$ws.Range("E46").Value = '  +12.45%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '141.23'
$ws.Range("D47").Style = "Normal"

$ws.Range("E47").Value = '  -0.94%  '

$ws.Range("E48").Value = '  -3.03%  '

$ws.Range("E49").Value = '  -2.54%  '

$ws.Range("B50").Value = 'InjectiveProtocol'

$ws.Range("C50").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '19.50'
$ws.Range("D50").Style = "Normal"

$ws.Range("E50").Value = '  -1.36%  '

$ws.Range("B51").Value = 'Hedera'

$ws.Range("C51").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0505'
$ws.Range("D51").Style = "Normal"

$ws.Range("E51").Value = '  -3.52%  '
